# Generate Report for Handback
#
# The localization-status report is regenerated: the e4599166-... file
# (row 3 of each table) has now been handed back and is in sync with
# en-US for both the zh-cn and de-de locales. Update the status cells,
# the "latest handback" timestamps, and clear the now-stale error detail
# message.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 3 is the e4599166-... file ---------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet: row 3 is the e4599166-... file -------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("K3").Value = "2016-08-19 16:52:29"
$wsZhCn.Range("P3").Value = ""
$wsZhCn.Columns.Item(16).ColumnWidth = 12.8

# --- de-de sheet: row 3 is the e4599166-... file -------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("K3").Value = "2016-08-19 16:52:36"
$wsDeDe.Range("P3").Value = ""
$wsDeDe.Columns.Item(16).ColumnWidth = 12.8
